$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "variable" field to "key" throughout the sheet:
#  - the header cell that literally says "variable"
#  - the "{variable}" placeholders embedded in other cells
# Update every cell referencing each shared string together so the
# string table entries are updated in place rather than duplicated.

# Column B (rows 2-4): "{variable}-codA" -> "{key}-codA"
foreach ($addr in @("B2", "B3", "B4")) {
    $cell = $ws.Range($addr)
    $text = [string]$cell.Value()
    $cell.Value = $text.Replace("variable", "key")
}

# Header cell C1: "variable" -> "key"
$ws.Range("C1").Value = "key"

# Column D (rows 2-4): "{variable} promotor fused with codA" -> "{key} promotor fused with codA"
foreach ($addr in @("D2", "D3", "D4")) {
    $cell = $ws.Range($addr)
    $text = [string]$cell.Value()
    $cell.Value = $text.Replace("variable", "key")
}

# Update the selection to match the target state (A4 selected)
$ws.Range("A4").Select()
